# Accommodation request upload and map with student
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 81: fill in the previously-blank trailer row with a real student record ---
$ws.Cells.Item(81, 1).Value2 = "CSIE"
$ws.Cells.Item(81, 2).Value2 = "Economic Informatics"
$ws.Cells.Item(81, 3).Value2 = "TC 419786"
$ws.Cells.Item(81, 4).NumberFormat = "0"
$ws.Cells.Item(81, 4).Value2 = 2971201360023
$ws.Cells.Item(81, 5).Value2 = "Eliza Ioana"
$ws.Cells.Item(81, 6).Value2 = "Țuțuianu"
$ws.Cells.Item(81, 7).Value2 = "L"
$ws.Cells.Item(81, 8).Value2 = "licenta_zi"
$ws.Cells.Item(81, 9).Value2 = 3
$ws.Cells.Item(81, 10).Value2 = $false
$ws.Cells.Item(81, 11).Value2 = $false
$ws.Cells.Item(81, 12).NumberFormat = "0.00"
$ws.Cells.Item(81, 12).Value2 = 10
$ws.Cells.Item(81, 13).Value2 = "F"
$ws.Cells.Item(81, 14).Value2 = "buget"
$ws.Cells.Item(81, 15).Value2 = 1076
$ws.Cells.Item(81, 16).Value2 = 180
$ws.Cells.Item(81, 17).NumberFormat = "0"
$ws.Cells.Item(81, 17).Value2 = 748318768

# --- Row 82: brand-new student record, including the new AccomodationRequestId ---
$ws.Cells.Item(82, 1).Value2 = "CSIE"
$ws.Cells.Item(82, 2).Value2 = "Economic Informatics"
$ws.Cells.Item(82, 3).Value2 = 456382
$ws.Cells.Item(82, 4).NumberFormat = "0"
$ws.Cells.Item(82, 4).Value2 = 2123456789123
$ws.Cells.Item(82, 5).Value2 = "Beatrice"
$ws.Cells.Item(82, 6).Value2 = "Vaduva"
$ws.Cells.Item(82, 7).Value2 = "D"
$ws.Cells.Item(82, 8).Value2 = "licenta_zi"
$ws.Cells.Item(82, 9).Value2 = 3
$ws.Cells.Item(82, 10).Value2 = $false
$ws.Cells.Item(82, 11).Value2 = $false
$ws.Cells.Item(82, 12).NumberFormat = "0.00"
$ws.Cells.Item(82, 12).Value2 = 10
$ws.Cells.Item(82, 13).Value2 = "F"
$ws.Cells.Item(82, 14).Value2 = "buget"
$ws.Cells.Item(82, 15).Value2 = 1076
$ws.Cells.Item(82, 16).Value2 = 180
$ws.Cells.Item(82, 17).NumberFormat = "0"
$ws.Cells.Item(82, 17).Value2 = 743456789

# --- New column R: "AccomodationRequestId" ---
$ws.Cells.Item(1, 18).Value2 = "AccomodationRequestId"

# Existing rows 2..80 get a sequential AccomodationRequestId (1..79)
for ($row = 2; $row -le 80; $row++) {
    $ws.Cells.Item($row, 18).Value2 = $row - 1
}

# Row 82 gets the last sequential id (80)
$ws.Cells.Item(82, 18).Value2 = 80

# --- Column widths for the new/affected columns (D and R), sized to fit their content ---
$ws.Columns.Item(4).ColumnWidth = 13.3
$ws.Columns.Item(17).ColumnWidth = 9.1

# --- Selection / scroll position to match the final saved view ---
$ws.Activate()
$ws.Range("C81").Select()
$excel.ActiveWindow.ScrollRow = 61
